$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = 5.074599999999998
    27  = 5.789000000000002
    32  = 6.589099999999995
    36  = 9.375600000000009
    38  = 4.766299999999995
    46  = 5.880000000000006
    54  = 4.963200000000004
    55  = 6.474599999999997
    56  = 4.617299999999997
    67  = 5.611499999999993
    69  = 5.283999999999993
    72  = 5.094800000000006
    83  = 5.481999999999998
    86  = 5.184500000000003
    91  = 4.969899999999996
    93  = 5.5226
    99  = 5.515699999999997
    104 = 9.427200000000001
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
